$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gearbox Tests")

# --- Defined names -------------------------------------------------------

# AIR_NUT -> AIR_NUT.point, and widen the range from a single cell to C9:E9
$n = $wb.Names.Item("AIR_NUT")
$n.RefersTo = '=''Gearbox Tests''!$C$9:$E$9'
$n.Name = 'AIR_NUT.point'

# DIPSTICK -> DIPSTICK.angle (range unchanged)
$n = $wb.Names.Item("DIPSTICK")
$n.Name = 'DIPSTICK.angle'

# New names for the HOUSING moments of inertia (centroidal), list-indexed
$wb.Names.Add('HOUSING.moments_of_inertia_centroidal.0', '=''Gearbox Tests''!$C$17')
$wb.Names.Add('HOUSING.moments_of_inertia_centroidal.1', '=''Gearbox Tests''!$C$18')
$wb.Names.Add('HOUSING.moments_of_inertia_centroidal.2', '=''Gearbox Tests''!$C$19')

# SHAFT_CENTERS keeps its existing name/range, but gains a ".distance" alias
# plus new per-axis point names
$wb.Names.Add('SHAFT_CENTERS.distance', '=''Gearbox Tests''!$C$10')
$wb.Names.Add('SHAFT_CENTERS.point_1.x', '=''Gearbox Tests''!$C$11')
$wb.Names.Add('SHAFT_CENTERS.point_1.y', '=''Gearbox Tests''!$C$12')
$wb.Names.Add('SHAFT_CENTERS.point_1.z', '=''Gearbox Tests''!$C$13')

# --- Worksheet data --------------------------------------------------------

# Row 9 (AIR_NUT.point): C9 becomes the first of a 3-point row
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 3

# Row 10 (SHAFT_CENTERS / SHAFT_CENTERS.distance): updated numeric value
$ws.Range("C10").Value = 159.99999999999989

# New rows 11-13: SHAFT_CENTERS.point_1.x/y/z
$ws.Range("B11").Value = "x"
$ws.Range("C11").Value = 899.99999999999909

$ws.Range("B12").Value = "y"
$ws.Range("C12").Value = 865.63897914069571

$ws.Range("B13").Value = "z"
$ws.Range("C13").Value = 261.00000000000006

# New rows 16-19: HOUSING MOI section
$ws.Range("B16").Value = "HOUSING MOI"

$ws.Range("B17").Value = "MXX"
$ws.Range("C17").Value = 2734036.863510197

$ws.Range("B18").Value = "MYY"
$ws.Range("C18").Value = 3833609.1842077454

$ws.Range("B19").Value = "MZZ"
$ws.Range("C19").Value = 4349785.7993760025

# --- View state: selection moves to D4 (matches the committed file) -------
[void]$ws.Range("D4").Select()
